# Apply "updated output files to reflect 11 holdouts" edit.
#
# Max Rose (NY-11) switched from NO to YES on impeachment. This requires:
#   1. Updating top_trump_dists (sheet7): mark Rose (row 8) as YES and give
#      him an announcement date of 2019-10-03 (serial 43741).
#   2. Removing Rose's row from full_list_of_nos (sheet12), shifting the
#      remaining holdout rows up.
#   3. Updating every summary/grouping worksheet's NO/YES counts (n) that
#      were derived from Rose's record flipping from NO to YES.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) prezresults2016 : R/NO count -1, R/YES count +1
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("prezresults2016")
$ws.Range("C3").Value = 10
$ws.Range("C5").Value = 21

# ---------------------------------------------------------------------
# 2) top_trump_dists : Rose (row 8) now voted YES, with an announce date
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("top_trump_dists")
$ws.Range("A8").Value = "YES"
$ws.Range("H8").Value = 43741
$ws.Range("H8").NumberFormat = "yyyy-mm-dd"
$ws.Range("J8").Value = 43741
$ws.Range("J8").NumberFormat = "yyyy-mm-dd"
$ws.Range("K8").Value = 10
$ws.Range("L8").Value = 2019

# ---------------------------------------------------------------------
# 3) full_list_of_nos : remove Rose's row (row 11); Torres Small and
#    Van Drew shift up from rows 12/13 to rows 11/12, sheet shrinks by 1
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("full_list_of_nos")
$ws.Rows.Item(11).Delete()

# ---------------------------------------------------------------------
# 4) gdp_vs_nationalavg : BELOW/NO count -1, BELOW/YES count +1
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("gdp_vs_nationalavg")
$ws.Range("C3").Value = 7
$ws.Range("C5").Value = 95

# ---------------------------------------------------------------------
# 5) college_vs_nationalavg : ABOVE/NO count -1, ABOVE/YES count +1
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("college_vs_nationalavg")
$ws.Range("C2").Value = 2
$ws.Range("C4").Value = 133

# ---------------------------------------------------------------------
# 6) nonwhite_vs_nationalavg : ABOVE/NO count -1, ABOVE/YES count +1
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("nonwhite_vs_nationalavg")
$ws.Range("C2").Value = 4
$ws.Range("C4").Value = 140

# ---------------------------------------------------------------------
# 7) rural_morethanfifth : BELOW/NO count -1, BELOW/YES count +1
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("rural_morethanfifth")
$ws.Range("C3").Value = 6
$ws.Range("C5").Value = 189

# ---------------------------------------------------------------------
# 8) margin_5_or_less : more_than_5_points/NO count -1, .../YES count +1
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("margin_5_or_less")
$ws.Range("C3").Value = 3
$ws.Range("C5").Value = 206
